$d = $word.ActiveDocument

$replacements = @(
    @{ old = "52×66=3432"; new = "75×74=5550" },
    @{ old = "44×37=1628"; new = "67×90=6030" },
    @{ old = "68×68=4624"; new = "71×55=3905" },
    @{ old = "84×20=1680"; new = "56×36=2016" },
    @{ old = "13×73=949";  new = "26×26=676"  },
    @{ old = "26×68=1768"; new = "65×31=2015" },
    @{ old = "71×60=4260"; new = "30×54=1620" },
    @{ old = "52×15=780";  new = "19×46=874"  },
    @{ old = "77×62=4774"; new = "63×74=4662" },
    @{ old = "30×13=390";  new = "83×79=6557" },
    @{ old = "41×56=2296"; new = "29×94=2726" },
    @{ old = "64×60=3840"; new = "21×83=1743" },
    @{ old = "49×57=2793"; new = "20×83=1660" },
    @{ old = "55×78=4290"; new = "68×44=2992" },
    @{ old = "72×35=2520"; new = "91×28=2548" },
    @{ old = "34×47=1598"; new = "70×42=2940" },
    @{ old = "27×27=729";  new = "42×29=1218" },
    @{ old = "63×66=4158"; new = "38×79=3002" },
    @{ old = "46×49=2254"; new = "37×66=2442" },
    @{ old = "37×48=1776"; new = "80×15=1200" },
    @{ old = "25×17=425";  new = "71×12=852"  },
    @{ old = "44×63=2772"; new = "48×13=624"  },
    @{ old = "97×13=1261"; new = "87×15=1305" },
    @{ old = "54×82=4428"; new = "51×63=3213" },
    @{ old = "17×99=1683"; new = "95×38=3610" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
